# Staging changes for main: add two new terrain-location columns (Desert /
# Badlands) to the lookup table on Sheet2, backed by six new shared strings
# (Desert, Oasis, Badlands, Cave, Spring, Abandoned Mine).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: Desert (filled top-to-bottom first, matches the original
#     authoring order so new shared-string indices land the same way) -----
$gValues = @(
    "Desert",
    "Tomb",
    "Fossil",
    "Small Village",
    "Ancient Battleground",
    "Monster Lair",
    "Ruined City",
    "Bandit Encampment",
    "Oasis"
)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = 1 + $i
    $cell = $ws.Range("G$row")
    $cell.Value = $gValues[$i]
    if ($row -gt 1) {
        $cell.HorizontalAlignment = -4131
    }
}

# --- Column H: Badlands (filled top-to-bottom after column G) ------------
$hValues = @(
    "Badlands",
    "Small Village",
    "Tomb",
    "Cave",
    "Ruined Fort",
    "Military Fort",
    "Monster Lair",
    "Bandit Encampment",
    "Ruined City",
    "Spring",
    "Abandoned Mine"
)
for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = 1 + $i
    $cell = $ws.Range("H$row")
    $cell.Value = $hValues[$i]
    if ($row -gt 1) {
        $cell.VerticalAlignment = -4160
    }
}

# --- Selection -------------------------------------------------------------
$ws.Range("C13").Select()
